$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 6.246
$ws.Range("D4").Value = -8.014999999999999

$ws.Range("D5").Value = -8.581999999999999

$ws.Range("B6").Value = 6.728
$ws.Range("D6").Value = -8.484999999999999

$ws.Range("B7").Value = 6.555

$ws.Range("B8").Value = 6.2
$ws.Range("D8").Value = -8.426

$ws.Range("B16").Value = 6.533999999999999
$ws.Range("D16").Value = -8.427000000000001

$ws.Range("B20").Value = 6.359

$ws.Range("B21").Value = 5.654999999999999

$ws.Range("D22").Value = -8.16
